$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the "Arveja Verde" market sheet.
# It is inserted as a new row 50, pushing the former rows 50-166 down to 51-167.
$ws.Rows("50:50").Insert()

$ws.Range("A50").Value = 5
$ws.Range("B50").Value = "Macroferia Regional de Talca"
$ws.Range("C50").Value = "Maule"
$ws.Range("D50").Value = 45260
$ws.Range("E50").Value = 7
$ws.Range("F50").Value = 100112022
$ws.Range("G50").Value = "Arveja Verde"
$ws.Range("H50").Value = "Sin especificar"
$ws.Range("I50").Value = "Primera"
$ws.Range("J50").Value = 500
$ws.Range("K50").Value = 18000
$ws.Range("L50").Value = 20000
$ws.Range("M50").Value = 18800
$ws.Range("N50").Value = "`$/saco 25 kilos"
$ws.Range("O50").Value = "Región del Maule"
$ws.Range("P50").Value = 752
$ws.Range("Q50").Value = 25
$ws.Range("R50").Value = "Hortaliza"
